$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 31.90834366666667
$ws.Range("H2").Value = 95.725031
$ws.Range("I2").Value = 0.1125536485145784
$ws.Range("J2").Value = 0.1157863270269485
$ws.Range("M2").Value = 14.349718
$ws.Range("N2").Value = 43.049154
$ws.Range("O2").Value = 0.1016415840981481
$ws.Range("P2").Value = 0.1034081666702025
$ws.Range("Q2").Value = 457.8757334637527
$ws.Range("R2").Value = 4120.881601173774
$ws.Range("S2").Value = 0.01144013113104792
$ws.Range("T2").Value = 0.01197325180333327
$ws.Range("G3").Value = 31.90834366666667
$ws.Range("H3").Value = 95.725031
$ws.Range("I3").Value = 0.1125536485145784
$ws.Range("J3").Value = 0.1157863270269485
$ws.Range("O3").Value = 0.04778708884009916
$ws.Range("P3").Value = 0.04861765281706964
$ws.Range("Q3").Value = 215.271618864474
$ws.Range("R3").Value = 1937.444569780266
$ws.Range("S3").Value = 0.005378611200843454
$ws.Range("T3").Value = 0.005629259448359871
$ws.Range("G4").Value = 31.90834366666667
$ws.Range("H4").Value = 95.725031
$ws.Range("I4").Value = 0.1125536485145784
$ws.Range("J4").Value = 0.1157863270269485
$ws.Range("M4").Value = 66.43651233333334
$ws.Range("N4").Value = 199.309537
$ws.Range("O4").Value = 0.4705815372480596
$ws.Range("P4").Value = 0.4787604843769264
$ws.Range("Q4").Value = 2119.879067546739
$ws.Range("R4").Value = 19078.91160792065
$ws.Range("S4").Value = 0.0529656689408681
$ws.Range("T4").Value = 0.05543391801164708
$ws.Range("G5").Value = 31.90834366666667
$ws.Range("H5").Value = 95.725031
$ws.Range("I5").Value = 0.1125536485145784
$ws.Range("J5").Value = 0.1157863270269485
$ws.Range("M5").Value = 7.2355625
$ws.Range("N5").Value = 14.471125
$ws.Range("O5").Value = 0.05125076564857627
$ws.Range("P5").Value = 0.03476102006337534
$ws.Range("Q5").Value = 230.8748148716458
$ws.Range("R5").Value = 1385.248889229875
$ws.Range("S5").Value = 0.005768460662912883
$ws.Range("T5").Value = 0.004024850836848296
$ws.Range("G6").Value = 31.90834366666667
$ws.Range("H6").Value = 95.725031
$ws.Range("I6").Value = 0.1125536485145784
$ws.Range("J6").Value = 0.1157863270269485
$ws.Range("M6").Value = 46.41124333333334
$ws.Range("N6").Value = 139.23373
$ws.Range("O6").Value = 0.328739024165117
$ws.Range("P6").Value = 0.3344526760724259
$ws.Range("Q6").Value = 1480.905902277292
$ws.Range("R6").Value = 13328.15312049563
$ws.Range("S6").Value = 0.03700077657890608
$ws.Range("T6").Value = 0.03872504692675999
$ws.Range("I7").Value = 0.2312918537506949
$ws.Range("J7").Value = 0.2379348388122522
$ws.Range("M7").Value = 14.349718
$ws.Range("N7").Value = 43.049154
$ws.Range("O7").Value = 0.1016415840981481
$ws.Range("P7").Value = 0.1034081666702025
$ws.Range("Q7").Value = 940.9106553002894
$ws.Range("R7").Value = 8468.195897702604
$ws.Range("S7").Value = 0.02350887040421781
$ws.Range("T7").Value = 0.02460440546854515
$ws.Range("I8").Value = 0.2312918537506949
$ws.Range("J8").Value = 0.2379348388122522
$ws.Range("O8").Value = 0.04778708884009916
$ws.Range("P8").Value = 0.04861765281706964
$ws.Range("S8").Value = 0.01105276436317568
$ws.Range("T8").Value = 0.01156783338645951
$ws.Range("I9").Value = 0.2312918537506949
$ws.Range("J9").Value = 0.2379348388122522
$ws.Range("M9").Value = 66.43651233333334
$ws.Range("N9").Value = 199.309537
$ws.Range("O9").Value = 0.4705815372480596
$ws.Range("P9").Value = 0.4787604843769264
$ws.Range("Q9").Value = 4356.240474929363
$ws.Range("R9").Value = 39206.16427436426
$ws.Range("S9").Value = 0.1088416760909554
$ws.Range("T9").Value = 0.1139137986798998
$ws.Range("I10").Value = 0.2312918537506949
$ws.Range("J10").Value = 0.2379348388122522
$ws.Range("M10").Value = 7.2355625
$ws.Range("N10").Value = 14.471125
$ws.Range("O10").Value = 0.05125076564857627
$ws.Range("P10").Value = 0.03476102006337534
$ws.Range("Q10").Value = 474.4356546477917
$ws.Range("R10").Value = 2846.61392788675
$ws.Range("S10").Value = 0.01185388459300164
$ws.Range("T10").Value = 0.008270857705728676
$ws.Range("I11").Value = 0.2312918537506949
$ws.Range("J11").Value = 0.2379348388122522
$ws.Range("M11").Value = 46.41124333333334
$ws.Range("N11").Value = 139.23373
$ws.Range("O11").Value = 0.328739024165117
$ws.Range("P11").Value = 0.3344526760724259
$ws.Range("Q11").Value = 3043.184080555998
$ws.Range("R11").Value = 27388.65672500398
$ws.Range("S11").Value = 0.07603465829934437
$ws.Range("T11").Value = 0.07957794357161907
$ws.Range("G12").Value = 85.57939900000001
$ws.Range("H12").Value = 256.738197
$ws.Range("I12").Value = 0.3018731932863474
$ws.Range("J12").Value = 0.3105433607867011
$ws.Range("M12").Value = 14.349718
$ws.Range("N12").Value = 43.049154
$ws.Range("O12").Value = 0.1016415840981481
$ws.Range("P12").Value = 0.1034081666702025
$ws.Range("Q12").Value = 1228.040242259482
$ws.Range("R12").Value = 11052.36218033534
$ws.Range("S12").Value = 0.03068286956239078
$ws.Range("T12").Value = 0.03211271961055602
$ws.Range("G13").Value = 85.57939900000001
$ws.Range("H13").Value = 256.738197
$ws.Range("I13").Value = 0.3018731932863474
$ws.Range("J13").Value = 0.3105433607867011
$ws.Range("O13").Value = 0.04778708884009916
$ws.Range("P13").Value = 0.04861765281706964
$ws.Range("Q13").Value = 577.3667212762381
$ws.Range("R13").Value = 5196.300491486142
$ws.Range("S13").Value = 0.01442564110601911
$ws.Range("T13").Value = 0.01509788929937383
$ws.Range("G14").Value = 85.57939900000001
$ws.Range("H14").Value = 256.738197
$ws.Range("I14").Value = 0.3018731932863474
$ws.Range("J14").Value = 0.3105433607867011
$ws.Range("M14").Value = 66.43651233333334
$ws.Range("N14").Value = 199.309537
$ws.Range("O14").Value = 0.4705815372480596
$ws.Range("P14").Value = 0.4787604843769264
$ws.Range("Q14").Value = 5685.596797142755
$ws.Range("R14").Value = 51170.3711742848
$ws.Range("S14").Value = 0.14205595135067
$ws.Range("T14").Value = 0.1486758898302796
$ws.Range("G15").Value = 85.57939900000001
$ws.Range("H15").Value = 256.738197
$ws.Range("I15").Value = 0.3018731932863474
$ws.Range("J15").Value = 0.3105433607867011
$ws.Range("M15").Value = 7.2355625
$ws.Range("N15").Value = 14.471125
$ws.Range("O15").Value = 0.05125076564857627
$ws.Range("P15").Value = 0.03476102006337534
$ws.Range("Q15").Value = 619.2150901769376
$ws.Range("R15").Value = 3715.290541061625
$ws.Range("S15").Value = 0.01547123228470595
$ws.Range("T15").Value = 0.01079480399485452
$ws.Range("G16").Value = 85.57939900000001
$ws.Range("H16").Value = 256.738197
$ws.Range("I16").Value = 0.3018731932863474
$ws.Range("J16").Value = 0.3105433607867011
$ws.Range("M16").Value = 46.41124333333334
$ws.Range("N16").Value = 139.23373
$ws.Range("O16").Value = 0.328739024165117
$ws.Range("P16").Value = 0.3344526760724259
$ws.Range("Q16").Value = 3971.846311309424
$ws.Range("R16").Value = 35746.61680178482
$ws.Range("S16").Value = 0.09923749898256158
$ws.Range("T16").Value = 0.103862058051637
$ws.Range("G17").Value = 23.7449455
$ws.Range("H17").Value = 47.489891
$ws.Range("I17").Value = 0.08375803763818537
$ws.Range("J17").Value = 0.05744244731349463
$ws.Range("M17").Value = 14.349718
$ws.Range("N17").Value = 43.049154
$ws.Range("O17").Value = 0.1016415840981481
$ws.Range("P17").Value = 0.1034081666702025
$ws.Range("Q17").Value = 340.733271850369
$ws.Range("R17").Value = 2044.399631102214
$ws.Range("S17").Value = 0.00851329962649747
$ws.Range("T17").Value = 0.00594001816573818
$ws.Range("G18").Value = 23.7449455
$ws.Range("H18").Value = 47.489891
$ws.Range("I18").Value = 0.08375803763818537
$ws.Range("J18").Value = 0.05744244731349463
$ws.Range("O18").Value = 0.04778708884009916
$ws.Range("P18").Value = 0.04861765281706964
$ws.Range("Q18").Value = 160.196747002371
$ws.Range("R18").Value = 961.180482014226
$ws.Range("S18").Value = 0.004002552785688333
$ws.Range("T18").Value = 0.002792716960450297
$ws.Range("G19").Value = 23.7449455
$ws.Range("H19").Value = 47.489891
$ws.Range("I19").Value = 0.08375803763818537
$ws.Range("J19").Value = 0.05744244731349463
$ws.Range("M19").Value = 66.43651233333334
$ws.Range("N19").Value = 199.309537
$ws.Range("O19").Value = 0.4705815372480596
$ws.Range("P19").Value = 0.4787604843769264
$ws.Range("Q19").Value = 1577.531364565078
$ws.Range("R19").Value = 9465.188187390468
$ws.Range("S19").Value = 0.03941498610865811
$ws.Range("T19").Value = 0.02750117389960476
$ws.Range("G20").Value = 23.7449455
$ws.Range("H20").Value = 47.489891
$ws.Range("I20").Value = 0.08375803763818537
$ws.Range("J20").Value = 0.05744244731349463
$ws.Range("M20").Value = 7.2355625
$ws.Range("N20").Value = 14.471125
$ws.Range("O20").Value = 0.05125076564857627
$ws.Range("P20").Value = 0.03476102006337534
$ws.Range("Q20").Value = 171.8080372243438
$ws.Range("R20").Value = 687.2321488973751
$ws.Range("S20").Value = 0.004292663558179269
$ws.Range("T20").Value = 0.001996758063553768
$ws.Range("G21").Value = 23.7449455
$ws.Range("H21").Value = 47.489891
$ws.Range("I21").Value = 0.08375803763818537
$ws.Range("J21").Value = 0.05744244731349463
$ws.Range("M21").Value = 46.41124333333334
$ws.Range("N21").Value = 139.23373
$ws.Range("O21").Value = 0.328739024165117
$ws.Range("P21").Value = 0.3344526760724259
$ws.Range("Q21").Value = 1102.032443537238
$ws.Range("R21").Value = 6612.194661223431
$ws.Range("S21").Value = 0.0275345355591622
$ws.Range("T21").Value = 0.01921178022414761
$ws.Range("G22").Value = 76.69186633333334
$ws.Range("H22").Value = 230.075599
$ws.Range("I22").Value = 0.270523266810194
$ws.Range("J22").Value = 0.2782930260606035
$ws.Range("M22").Value = 14.349718
$ws.Range("N22").Value = 43.049154
$ws.Range("O22").Value = 0.1016415840981481
$ws.Range("P22").Value = 0.1034081666702025
$ws.Range("Q22").Value = 1100.506654777028
$ws.Range("R22").Value = 9904.559892993248
$ws.Range("S22").Value = 0.02749641337399408
$ws.Range("T22").Value = 0.0287777716220299
$ws.Range("G23").Value = 76.69186633333334
$ws.Range("H23").Value = 230.075599
$ws.Range("I23").Value = 0.270523266810194
$ws.Range("J23").Value = 0.2782930260606035
$ws.Range("O23").Value = 0.04778708884009916
$ws.Range("P23").Value = 0.04861765281706964
$ws.Range("Q23").Value = 517.406431113546
$ws.Range("R23").Value = 4656.657880021914
$ws.Range("S23").Value = 0.01292751938437259
$ws.Range("T23").Value = 0.01352995372242614
$ws.Range("G24").Value = 76.69186633333334
$ws.Range("H24").Value = 230.075599
$ws.Range("I24").Value = 0.270523266810194
$ws.Range("J24").Value = 0.2782930260606035
$ws.Range("M24").Value = 66.43651233333334
$ws.Range("N24").Value = 199.309537
$ws.Range("O24").Value = 0.4705815372480596
$ws.Range("P24").Value = 0.4787604843769264
$ws.Range("Q24").Value = 5095.140123520852
$ws.Range("R24").Value = 45856.26111168767
$ws.Range("S24").Value = 0.1273032547569081
$ws.Range("T24").Value = 0.1332357039554951
$ws.Range("G25").Value = 76.69186633333334
$ws.Range("H25").Value = 230.075599
$ws.Range("I25").Value = 0.270523266810194
$ws.Range("J25").Value = 0.2782930260606035
$ws.Range("M25").Value = 7.2355625
$ws.Range("N25").Value = 14.471125
$ws.Range("O25").Value = 0.05125076564857627
$ws.Range("P25").Value = 0.03476102006337534
$ws.Range("Q25").Value = 554.9087920964793
$ws.Range("R25").Value = 3329.452752578875
$ws.Range("S25").Value = 0.01386452454977652
$ws.Range("T25").Value = 0.009673749462390075
$ws.Range("G26").Value = 76.69186633333334
$ws.Range("H26").Value = 230.075599
$ws.Range("I26").Value = 0.270523266810194
$ws.Range("J26").Value = 0.2782930260606035
$ws.Range("M26").Value = 46.41124333333334
$ws.Range("N26").Value = 139.23373
$ws.Range("O26").Value = 0.328739024165117
$ws.Range("P26").Value = 0.3344526760724259
$ws.Range("Q26").Value = 3559.364870083808
$ws.Range("R26").Value = 32034.28383075427
$ws.Range("S26").Value = 0.08893155474514275
$ws.Range("T26").Value = 0.09307584729826221
